$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Demo1 scene): CamOffestPos / CamOffestRot
$ws.Range("J2").Value = "0,8,7"
$ws.Range("K2").Value = "45,180"

# Row 6 (City / SelectScene): CamOffestPos / CamOffestRot
$ws.Range("J6").Value = "0,8,-7"
$ws.Range("K6").Value = "45,0"

# Update the active selection to match the saved view state
$ws.Range("K7").Select()
